$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B9 currently holds "Plan de Pruebas" (E108 row); change it to the new
# configuration-management plan text.
$ws.Range("B9").Value = "Plan de Gestión de la configuración"

# Cell B10 (E109 row) was empty; now holds the "Plan de Pruebas" text that used
# to live in B9.
$ws.Range("B10").Value = "Plan de Pruebas"

# Update the active selection to match the edited workbook.
$ws.Range("F11").Select()
